# Updated cryptos list on Thu Jan 18 23:42:02 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest snapshot.
# Price cells are forced to Text (NumberFormat "@") before the write so that
# numeric-looking strings (e.g. "313.97", "9.70") are stored as literal text
# instead of being auto-coerced into floating point numbers (which would
# silently drop meaningful trailing zeros). ClearFormats() afterwards removes
# the temporary Text number-format so the cell keeps its original (default)
# style - only the stored value/type changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "41.317.19"
$c.ClearFormats()
$ws.Range("E2").Value = "  -3.29%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.466.90"
$c.ClearFormats()
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "313.97"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.45%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "94.41"
$c.ClearFormats()
$ws.Range("E6").Value = "  -6.96%  "
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("E8").Value = "  +0.08%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.502"
$c.ClearFormats()
$ws.Range("E9").Value = "  -4.56%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "33.58"
$c.ClearFormats()
$ws.Range("E10").Value = "  -6.67%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  -0.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.02"
$c.ClearFormats()
$ws.Range("E13").Value = "  -4.00%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.850.51"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.10%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.469.78"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.10%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.68"
$c.ClearFormats()
$ws.Range("E16").Value = "  -6.21%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.786"
$c.ClearFormats()
$ws.Range("E17").Value = "  -2.53%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "41.307.53"
$c.ClearFormats()
$ws.Range("E18").Value = "  -3.24%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.ClearFormats()
$ws.Range("E19").Value = "  -6.24%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0920"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.09%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.48"
$c.ClearFormats()
$ws.Range("E21").Value = "  -5.44%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.05"
$c.ClearFormats()
$ws.Range("E22").Value = "  -1.78%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "237.57"
$c.ClearFormats()
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("E25").Value = "  -5.15%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.43"
$c.ClearFormats()
$ws.Range("E27").Value = "  -6.46%  "
$ws.Range("E28").Value = "  -4.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.70"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.31%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "36.02"
$c.ClearFormats()
$ws.Range("E30").Value = "  -7.70%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "153.09"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.44%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.58"
$c.ClearFormats()
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("E33").Value = "  -6.28%  "
$ws.Range("E34").Value = "  -0.89%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0754"
$c.ClearFormats()
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("E36").Value = "  -4.97%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.ClearFormats()
$ws.Range("E37").Value = "  -6.89%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "16.89"
$c.ClearFormats()
$ws.Range("E38").Value = "  -7.22%  "
$ws.Range("E39").Value = "  -6.46%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.27"
$c.ClearFormats()
$ws.Range("E41").Value = "  +1.83%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "21.31"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("E43").Value = "  +0.23%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.991.62"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.06%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0285"
$c.ClearFormats()
$ws.Range("E45").Value = "  -4.64%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.07"
$c.ClearFormats()
$ws.Range("E46").Value = "  -6.37%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.74"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.30%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "69.81"
$c.ClearFormats()
$ws.Range("E48").Value = "  -3.06%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "75.91"
$c.ClearFormats()
$ws.Range("E49").Value = "  -5.37%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "96.92"
$c.ClearFormats()
$ws.Range("E50").Value = "  -4.18%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.179"
$c.ClearFormats()
$ws.Range("E51").Value = "  -6.08%  "
